# Tried to implement Penalty Reward System (unfinished)
# - Adjust a handful of weekly/monthly "Requested quantity" figures
# - Remove a few rows that no longer apply (penalty/reward rows retired)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Weekly Quantity"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Update quantities first (while row numbers still match the original layout)
$ws1.Range("B8").Value = 40     # was 230
$ws1.Range("B9").Value = 190    # was 270
$ws1.Range("B13").Value = 560   # was 710
$ws1.Range("B15").Value = 450   # was 830

# Remove rows that are no longer needed, bottom-to-top so earlier row
# numbers stay valid while we work.
$ws1.Range("A34:A36").EntireRow.Delete()   # weeks of 45361/45368/45375
$ws1.Rows(14).Delete()                      # week of 45130
$ws1.Rows(12).Delete()                      # week of 45116

# ---------------------------------------------------------------------
# Sheet 2: "Monthly Trend"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws2.Range("B4").Value = 970    # was 1240
$ws2.Range("B5").Value = 1010   # was 2640

$ws2.Rows(13).Delete()           # month of 45382
